$wb = $excel.ActiveWorkbook

# --- Arkusz1: add the new "Boulangerie & Pâtisserie" subcategories in column F ---
$ws1 = $wb.Worksheets.Item("Arkusz1")
$ws1.Range("F2").Value = "Pain"
$ws1.Range("F3").Value = "Pains spéciaux & précuits"
$ws1.Range("F4").Value = "Viennoiseries"

# --- Arkusz2: move the cursor/selection to B14 ---
$ws2 = $wb.Worksheets.Item("Arkusz2")
[void]$ws2.Range("B14").Select()

# --- Arkusz1: move the cursor/selection to F9 (also re-activates Arkusz1 as the visible tab) ---
[void]$ws1.Range("F9").Select()
